# #327 Ajout des profils d'acces a58d18c1e8091c98efec92c8c093b361a253eee5
#
# Changes applied:
#  1. Metadata sheet: bump the "Date" value (row 8, column B).
#  2. Elements sheet: the "Mapping: RIM Mapping" column (AK / 37) and the
#     "Mapping: Spécification métier vers l'extension ROR OrganizationPrice"
#     column (AL / 38) are swapped for every row (header + data), because the
#     sheet author re-ordered/re-purposed those two trailing mapping columns.
#  3. Row 39 (Extension.extension:othersDeliveryIncluded.value[x]) gets a
#     corrected business mapping instead of the stale duplicated text that
#     row 34 also has, so it is special-cased after the generic swap.
#  4. The two mapping columns swap their (approximate) display widths to
#     match their new, swapped content.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata!B8 - Date
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value2 = "2024-03-19T13:17:15+00:00"

# ---------------------------------------------------------------------------
# 2. Elements - swap columns AK (37, "Mapping: RIM Mapping") and
#    AL (38, "Mapping: Spécification métier vers l'extension ROR
#    OrganizationPrice") for every row, including the header row.
# ---------------------------------------------------------------------------
$els = $wb.Worksheets.Item("Elements")

for ($r = 1; $r -le 86; $r++) {
    $akCell = $els.Cells.Item($r, 37)
    $alCell = $els.Cells.Item($r, 38)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    $akCell.Value2 = $alVal
    $alCell.Value2 = $akVal
}

# ---------------------------------------------------------------------------
# 3. Row 39 correction - the business mapping is no longer the stale
#    "...prestationsNonObligatoiresIncluses" text shared with row 34, it
#    becomes a dedicated value, and the RIM mapping goes back to "N/A".
# ---------------------------------------------------------------------------
$els.Cells.Item(39, 37).Value2 = "ForfaitSocleHebergement.autresPrestationsNonObligatoiresIncluses"
$els.Cells.Item(39, 38).Value2 = "N/A"

# ---------------------------------------------------------------------------
# 4. Swap the two columns' widths to mirror their new content
#    (narrow "N/A"-style RIM column vs. wide free-text business column).
# ---------------------------------------------------------------------------
$els.Columns.Item(37).ColumnWidth = 133.0
$els.Columns.Item(38).ColumnWidth = 24.17
